{"js": "// Replace the 25 two-digit multiplication problems in the table cells.\n// Each (old, new) pair is unique, so a simple search+replace per pair is unambiguous.\nconst replacements = [\n  [\"14\u00d763=\", \"17\u00d751=\"],\n  [\"24\u00d761=\", \"74\u00d749=\"],\n  [\"77\u00d738=\", \"68\u00d725=\"],\n  [\"65\u00d747=\", \"55\u00d745=\"],\n  [\"72\u00d713=\", \"63\u00d756=\"],\n  [\"23\u00d795=\", \"31\u00d756=\"],\n  [\"71\u00d753=\", \"13\u00d748=\"],\n  [\"31\u00d721=\", \"88\u00d779=\"],\n  [\"97\u00d728=\", \"98\u00d717=\"],\n  [\"82\u00d751=\", \"55\u00d740=\"],\n  [\"17\u00d754=\", \"47\u00d783=\"],\n  [\"24\u00d795=\", \"75\u00d756=\"],\n  [\"98\u00d792=\", \"92\u00d760=\"],\n  [\"70\u00d771=\", \"44\u00d745=\"],\n  [\"96\u00d782=\", \"97\u00d725=\"],\n  [\"78\u00d726=\", \"68\u00d751=\"],\n  [\"95\u00d784=\", \"45\u00d749=\"],\n  [\"85\u00d768=\", \"37\u00d798=\"],\n  [\"68\u00d776=\", \"44\u00d779=\"],\n  [\"93\u00d749=\", \"89\u00d739=\"],\n  [\"78\u00d751=\", \"70\u00d773=\"],\n  [\"27\u00d730=\", \"75\u00d724=\"],\n  [\"13\u00d756=\", \"39\u00d744=\"],\n  [\"11\u00d742=\", \"21\u00d712=\"],\n  [\"14\u00d793=\", \"99\u00d776=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 two-digit multiplication problems in the table cells.\n# Each (old, new) pair is unique in the document, so Find/Replace-All per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @{Old = \"14\u00d763=\"; New = \"17\u00d751=\"}\n    @{Old = \"24\u00d761=\"; New = \"74\u00d749=\"}\n    @{Old = \"77\u00d738=\"; New = \"68\u00d725=\"}\n    @{Old = \"65\u00d747=\"; New = \"55\u00d745=\"}\n    @{Old = \"72\u00d713=\"; New = \"63\u00d756=\"}\n    @{Old = \"23\u00d795=\"; New = \"31\u00d756=\"}\n    @{Old = \"71\u00d753=\"; New = \"13\u00d748=\"}\n    @{Old = \"31\u00d721=\"; New = \"88\u00d779=\"}\n    @{Old = \"97\u00d728=\"; New = \"98\u00d717=\"}\n    @{Old = \"82\u00d751=\"; New = \"55\u00d740=\"}\n    @{Old = \"17\u00d754=\"; New = \"47\u00d783=\"}\n    @{Old = \"24\u00d795=\"; New = \"75\u00d756=\"}\n    @{Old = \"98\u00d792=\"; New = \"92\u00d760=\"}\n    @{Old = \"70\u00d771=\"; New = \"44\u00d745=\"}\n    @{Old = \"96\u00d782=\"; New = \"97\u00d725=\"}\n    @{Old = \"78\u00d726=\"; New = \"68\u00d751=\"}\n    @{Old = \"95\u00d784=\"; New = \"45\u00d749=\"}\n    @{Old = \"85\u00d768=\"; New = \"37\u00d798=\"}\n    @{Old = \"68\u00d776=\"; New = \"44\u00d779=\"}\n    @{Old = \"93\u00d749=\"; New = \"89\u00d739=\"}\n    @{Old = \"78\u00d751=\"; New = \"70\u00d773=\"}\n    @{Old = \"27\u00d730=\"; New = \"75\u00d724=\"}\n    @{Old = \"13\u00d756=\"; New = \"39\u00d744=\"}\n    @{Old = \"11\u00d742=\"; New = \"21\u00d712=\"}\n    @{Old = \"14\u00d793=\"; New = \"99\u00d776=\"}\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair.New, $wdReplaceAll)\n}\n\n"}
